$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stale test data (rows 2-7) while preserving the header row
# and the numeric-index formatting already applied to column A.
$ws.Range("A2:E7").ClearContents()

# Extend the existing "index" style/formatting from A7 down through the
# two new rows that will be added (rows 8 and 9).
$ws.Range("A7").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Column A: running index 0..7 ---
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7

# --- Column B: NAME ---
$ws.Range("B2").Value = "Test1"
$ws.Range("B3").Value = "Test2"
$ws.Range("B4").Value = "Test3"
$ws.Range("B5").Value = "Test4"
$ws.Range("B6").Value = "Test5"
$ws.Range("B7").Value = "Test6"
$ws.Range("B8").Value = "Test7"
$ws.Range("B9").Value = "Test8"

# --- Column C: SENTENCES id ---
$ws.Range("C2").Value = "P1_W1_S4"
$ws.Range("C3").Value = "P1_W1_S3"
$ws.Range("C4").Value = "P1_W1_S2"
$ws.Range("C5").Value = "P1_W1_S1"
$ws.Range("C6").Value = "P1_W2_S4"
$ws.Range("C7").Value = "P1_W2_S3"
$ws.Range("C8").Value = "P1_W2_S2"
$ws.Range("C9").Value = "P1_W2_S1"

# --- Column D: INPUT_SENTENCE text ---
$ws.Range("D2").Value = "Enjoy the fair weather while in the tropics."
$ws.Range("D3").Value = 'You''re used to being on the field.'
$ws.Range("D4").Value = "The ballet is about to begin."
$ws.Range("D5").Value = "We picked grapes for wine"
$ws.Range("D6").Value = "he is capable and willing to make decisions."
$ws.Range("D7").Value = "Big muscles are not necessarily strong ones"
$ws.Range("D8").Value = "You want him to do well"
$ws.Range("D9").Value = 'I think I''m getting better.'

# --- Column E: INTELLIGIBILITY_SCORE ---
$ws.Range("E2").Value = 0.081632653061224483
$ws.Range("E3").Value = 0.15384615384615391
$ws.Range("E4").Value = 0.1764705882352941
$ws.Range("E5").Value = 0.1333333333333333
$ws.Range("E6").Value = 0.1224489795918367
$ws.Range("E7").Value = 0.125
$ws.Range("E8").Value = 0.14285714285714279
$ws.Range("E9").Value = 0.125

$ws.Range("A1:E9").EntireColumn.AutoFit()
$ws.Range("C7").Select()

$wb.Save()
